$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'286.92"
$ws.Range("E2").Value = "'1.38%"
$ws.Range("G2").Value = "'12"

# Row 3
$ws.Range("D3").Value = "'29.40"
$ws.Range("E3").Value = "'3.89%"
$ws.Range("G3").Value = "'12"

# Row 4
$ws.Range("D4").Value = "'5.068"
$ws.Range("E4").Value = "'0.89%"
$ws.Range("G4").Value = "'12"

# Row 5
$ws.Range("D5").Value = "'0.06733"
$ws.Range("E5").Value = "'3.31%"
$ws.Range("G5").Value = "'12"

# Row 6
$ws.Range("D6").Value = "'7.334"
$ws.Range("E6").Value = "'1.32%"
$ws.Range("G6").Value = "'12"

# Row 7
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.381"
$ws.Range("E7").Value = "'-0.72%"
$ws.Range("G7").Value = "'12"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9049"
$ws.Range("E8").Value = "'-0.74%"
$ws.Range("G8").Value = "'12"

# Row 9
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1590"
$ws.Range("E9").Value = "'2.03%"
$ws.Range("G9").Value = "'12"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.06962"
$ws.Range("E10").Value = "'6.53%"
$ws.Range("G10").Value = "'12"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07642"
$ws.Range("E11").Value = "'0.98%"
$ws.Range("G11").Value = "'12"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02925"
$ws.Range("E12").Value = "'6.18%"
$ws.Range("G12").Value = "'12"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.08994"
$ws.Range("E13").Value = "'0.29%"
$ws.Range("G13").Value = "'12"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001577"
$ws.Range("E14").Value = "'-1.49%"
$ws.Range("G14").Value = "'12"

# Row 15
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04480"
$ws.Range("E15").Value = "'1.70%"
$ws.Range("G15").Value = "'12"

# Row 16
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0006467"
$ws.Range("E16").Value = "'1.54%"
$ws.Range("G16").Value = "'12"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006260"
$ws.Range("E17").Value = "'3.57%"
$ws.Range("G17").Value = "'12"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.450"
$ws.Range("E18").Value = "'0.24%"
$ws.Range("G18").Value = "'12"

# Row 19
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.449"
$ws.Range("E19").Value = "'2.64%"
$ws.Range("G19").Value = "'12"

# Row 20
$ws.Range("D20").Value = "'2.228"
$ws.Range("E20").Value = "'-0.67%"
$ws.Range("G20").Value = "'12"

# Row 21
$ws.Range("D21").Value = "'0.3205"
$ws.Range("E21").Value = "'0.59%"
$ws.Range("G21").Value = "'12"

# Row 22
$ws.Range("D22").Value = "'0.1315"
$ws.Range("E22").Value = "'2.60%"
$ws.Range("G22").Value = "'12"

# Row 23
$ws.Range("D23").Value = "'4.046"
$ws.Range("E23").Value = "'2.04%"
$ws.Range("G23").Value = "'12"

# Row 24
$ws.Range("D24").Value = "'0.1581"
$ws.Range("E24").Value = "'2.43%"
$ws.Range("G24").Value = "'12"

# Row 25
$ws.Range("D25").Value = "'0.001196"
$ws.Range("E25").Value = "'1.05%"
$ws.Range("G25").Value = "'12"

# Row 26
$ws.Range("D26").Value = "'0.004377"
$ws.Range("E26").Value = "'-1.08%"
$ws.Range("G26").Value = "'12"

# Row 27
$ws.Range("D27").Value = "'0.0001201"
$ws.Range("E27").Value = "'0.11%"
$ws.Range("G27").Value = "'12"

# Row 28
$ws.Range("D28").Value = "'0.0001617"
$ws.Range("E28").Value = "'-0.05%"
$ws.Range("G28").Value = "'12"

# Row 29
$ws.Range("G29").Value = "'12"

# Row 30
$ws.Range("G30").Value = "'12"

# Row 31
$ws.Range("G31").Value = "'12"

# Row 32
$ws.Range("G32").Value = "'12"

# Row 33
$ws.Range("G33").Value = "'12"

# Row 34
$ws.Range("G34").Value = "'12"

# Row 35
$ws.Range("G35").Value = "'12"

# Row 36
$ws.Range("G36").Value = "'12"

# Row 37
$ws.Range("G37").Value = "'12"

# Row 38
$ws.Range("G38").Value = "'12"

# Row 39
$ws.Range("G39").Value = "'12"

# Row 40
$ws.Range("D40").Value = "'0.04254"
$ws.Range("E40").Value = "'3.07%"
$ws.Range("G40").Value = "'12"

# Row 41
$ws.Range("D41").Value = "'0.006810"
$ws.Range("E41").Value = "'2.29%"
$ws.Range("G41").Value = "'12"

# Row 42
$ws.Range("D42").Value = "'0.1240"
$ws.Range("E42").Value = "'1.00%"
$ws.Range("G42").Value = "'12"

# Row 43
$ws.Range("D43").Value = "'0.002231"
$ws.Range("E43").Value = "'2.41%"
$ws.Range("G43").Value = "'12"

# Row 44
$ws.Range("D44").Value = "'0.01290"
$ws.Range("E44").Value = "'6.84%"
$ws.Range("G44").Value = "'12"

# Row 45
$ws.Range("D45").Value = "'0.00005696"
$ws.Range("E45").Value = "'1.24%"
$ws.Range("G45").Value = "'12"

# Row 46
$ws.Range("D46").Value = "'1.963"
$ws.Range("E46").Value = "'-0.16%"
$ws.Range("G46").Value = "'12"

# Row 47
$ws.Range("D47").Value = "'0.01499"
$ws.Range("E47").Value = "'-19.00%"
$ws.Range("G47").Value = "'12"

# Row 48
$ws.Range("G48").Value = "'12"

# Row 49
$ws.Range("G49").Value = "'12"

# Row 50
$ws.Range("G50").Value = "'12"

# Row 51
$ws.Range("G51").Value = "'12"
